# Generate Report for Handback
#
# Applies the "handback" report-generation edit:
#  - "In Translation" status text is replaced everywhere by
#    "Handed back: in sync with en-US" (Overview!E2/F2/E3/F3 and the
#    Status column (C) on the zh-cn / de-de report sheets).
#  - The zh-cn / de-de report sheets get their "Latest Target File" (I)
#    and "Latest Handback File" (J) columns populated for both data rows,
#    with I turned into a hyperlink to the source .md file (mirroring the
#    existing Source File Name hyperlink style/column A).
#  - The "Latest Handback DateTime" (K) column is stamped: zh-cn gets
#    2016-09-02 14:29:54, de-de gets 2016-09-02 14:30:05.
#  - A few columns are widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$hyperlinkFontColor = 15570276   # OLE BGR for RGB(0x64,0x95,0xED) - matches the workbook's existing HyperLink style

function Set-ColWidthApprox($ws, [int]$colIndex, [double]$targetStoredWidth) {
    # The host's internal column-width unit is quantized to 1/6ths once it
    # round-trips through ColumnWidth; pick the COM input that lands on the
    # stored value closest to the width used by the reference workbook.
    $n = [Math]::Round($targetStoredWidth * 6 - 5)
    $cw = $n / 6.0
    $ws.Columns.Item($colIndex).ColumnWidth = $cw
}

function Set-HandbackHyperlink($ws, [string]$cellRef, [string]$fileName) {
    $target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bcdd7db6bfb3d300e4895cb8fe502d939ffd8d88/e2e/" + $fileName
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $fileName)
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkFontColor
}

# ---------------------------------------------------------------------
# Overview sheet: "In Translation" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

Set-ColWidthApprox $overview 5 29.9777047293527
Set-ColWidthApprox $overview 6 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn report sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
Set-ColWidthApprox $zhcn 3 29.9777047293527

Set-HandbackHyperlink $zhcn "I2" "53104752-f2d3-429b-aad1-20993a0a468d.md"
Set-HandbackHyperlink $zhcn "I3" "b4472403-a166-40a5-b7f0-03a53221021c.md"

$zhcn.Range("J2").Value = "53104752-f2d3-429b-aad1-20993a0a468d.e31c8bafbf178f1c48f741cde553945f7c3cb11d.zh-cn.xlf"
$zhcn.Range("J3").Value = "b4472403-a166-40a5-b7f0-03a53221021c.067966c5a7e6d2afe61f8434292f1ec282feedc8.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-09-02 14:29:54"
$zhcn.Range("K3").Value = "2016-09-02 14:29:54"

Set-ColWidthApprox $zhcn 9 40
Set-ColWidthApprox $zhcn 10 40

# ---------------------------------------------------------------------
# de-de report sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"
Set-ColWidthApprox $dede 3 29.9777047293527

Set-HandbackHyperlink $dede "I2" "53104752-f2d3-429b-aad1-20993a0a468d.md"
Set-HandbackHyperlink $dede "I3" "b4472403-a166-40a5-b7f0-03a53221021c.md"

$dede.Range("J2").Value = "53104752-f2d3-429b-aad1-20993a0a468d.e31c8bafbf178f1c48f741cde553945f7c3cb11d.de-de.xlf"
$dede.Range("J3").Value = "b4472403-a166-40a5-b7f0-03a53221021c.067966c5a7e6d2afe61f8434292f1ec282feedc8.de-de.xlf"

$dede.Range("K2").Value = "2016-09-02 14:30:05"
$dede.Range("K3").Value = "2016-09-02 14:30:05"

Set-ColWidthApprox $dede 9 40
Set-ColWidthApprox $dede 10 40
